$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date column C for all existing data rows (2-487)
#    from 45189 (2023-09-20) to 45190 (2023-09-21)
$ws.Range("C2:C487").Value = 45190

# 2. Append new row 488
$ws.Range("A488").Value = "A 44372-2023"
$ws.Range("B488").Value = 45188
$ws.Range("C488").Value = 45190
$ws.Range("D488").Value = "VÄSTERBOTTENS LÄN"
$ws.Range("E488").Value = "ROBERTSFORS"
$ws.Range("G488").Value = 5.4
$ws.Range("H488").Value = 0
$ws.Range("I488").Value = 0
$ws.Range("J488").Value = 0
$ws.Range("K488").Value = 0
$ws.Range("L488").Value = 0
$ws.Range("M488").Value = 0
$ws.Range("N488").Value = 0
$ws.Range("O488").Value = 0
$ws.Range("P488").Value = 0
$ws.Range("Q488").Value = 0
$ws.Range("R488").WrapText = $true

# 3. Append new row 489
$ws.Range("A489").Value = "A 44374-2023"
$ws.Range("B489").Value = 45188
$ws.Range("C489").Value = 45190
$ws.Range("D489").Value = "VÄSTERBOTTENS LÄN"
$ws.Range("E489").Value = "ROBERTSFORS"
$ws.Range("G489").Value = 2.9
$ws.Range("H489").Value = 0
$ws.Range("I489").Value = 0
$ws.Range("J489").Value = 0
$ws.Range("K489").Value = 0
$ws.Range("L489").Value = 0
$ws.Range("M489").Value = 0
$ws.Range("N489").Value = 0
$ws.Range("O489").Value = 0
$ws.Range("P489").Value = 0
$ws.Range("Q489").Value = 0
$ws.Range("R489").WrapText = $true

# 4. Apply the same date number format used by the rest of column B / C to the new cells
$ws.Range("B488:C489").NumberFormat = "YYYY-MM-DD"

# 5. Row heights: row 487 and row 488 get an explicit 15pt custom height (row 489, the
#    new last row, stays without an explicit height - matching the target workbook)
$ws.Range("A487:A488").EntireRow.RowHeight = 15
